$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for de1badf3... row
$wsOverview.Range("G2").Value = "2016-09-05 21:15:34"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime for de1badf3... row
$wsZhCn.Range("H2").Value = "2016-09-05 21:15:29"
$wsZhCn.Range("K2").Value = "2016-09-05 21:15:46"

# de-de sheet: Correspond Handoff Datetime (shared with Overview) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-09-05 21:15:34"
$wsDeDe.Range("K2").Value = "2016-09-05 21:15:54"
